$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 7
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 7
$ws.Range("F18").Value = 5
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 3
$ws.Range("F23").Value = 7
$ws.Range("F24").Value = 5
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 4
$ws.Range("F28").Value = -4
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 3
$ws.Range("F32").Value = -1
$ws.Range("F33").Value = -2
$ws.Range("F35").Value = 3
